# Scheduled runner update: refresh cached market-board derived profit figures
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 15102.857
$ws.Range("I34").Value = 10953.5
$ws.Range("K34").Value = 10953.5
$ws.Range("M34").Value = -10750.5
$ws.Range("H36").Value = 15102.857
$ws.Range("I36").Value = 10953.5
$ws.Range("K36").Value = 10953.5
$ws.Range("M36").Value = -10238.5
$ws.Range("H40").Value = 5225
$ws.Range("I40").Value = 4667
$ws.Range("K40").Value = 4667
$ws.Range("M40").Value = -4492
$ws.Range("H62").Value = 15486.258
$ws.Range("I62").Value = 22171.285
$ws.Range("K62").Value = 22171.285
$ws.Range("M62").Value = -21547.285
$ws.Range("H65").Value = 15486.258
$ws.Range("I65").Value = 22171.285
$ws.Range("K65").Value = 110856.425
$ws.Range("M65").Value = -107736.425
$ws.Range("H70").Value = 1977.2142
$ws.Range("J70").Value = 2006.3077
$ws.Range("L70").Value = 6018.9231
$ws.Range("N70").Value = -6558.9231
$ws.Range("H73").Value = 1977.2142
$ws.Range("J73").Value = 2006.3077
$ws.Range("L73").Value = 6018.9231
$ws.Range("N73").Value = -7890.9231
$ws.Range("H76").Value = 8273.857
$ws.Range("I76").Value = 8484.5
$ws.Range("K76").Value = 8484.5
$ws.Range("M76").Value = -8169.5
$ws.Range("H79").Value = 8273.857
$ws.Range("I79").Value = 8484.5
$ws.Range("K79").Value = 8484.5
$ws.Range("M79").Value = -7392.5
$ws.Range("H116").Value = 3000005
$ws.Range("I116").Value = 3000005
$ws.Range("K116").Value = 3000005
$ws.Range("M116").Value = -2996563
$ws.Range("H132").Value = 3565.3872
$ws.Range("I132").Value = 2657.625
$ws.Range("K132").Value = 7972.875
$ws.Range("M132").Value = -5442.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11709.901
$ws.Range("I32").Value = 8740.223
$ws.Range("K32").Value = 8740.223
$ws.Range("M32").Value = -8453.223
$ws.Range("H61").Value = 3922.1785
$ws.Range("I61").Value = 3743.236
$ws.Range("K61").Value = 3743.236
$ws.Range("M61").Value = -3531.236
$ws.Range("H136").Value = 3922.1785
$ws.Range("I136").Value = 3743.236
$ws.Range("K136").Value = 11229.708
$ws.Range("M136").Value = -8679.707999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 7185.3335
$ws.Range("I19").Value = 773
$ws.Range("K19").Value = 773
$ws.Range("M19").Value = -600
$ws.Range("H99").Value = 6250
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 6250
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 6250
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -9246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 89.16
$ws.Range("I7").Value = 77.5
$ws.Range("J7").Value = 104
$ws.Range("K7").Value = 77.5
$ws.Range("L7").Value = 104
$ws.Range("M7").Value = 35.5
$ws.Range("N7").Value = -330
$ws.Range("H16").Value = 999.2
$ws.Range("I16").Value = 888
$ws.Range("K16").Value = 888
$ws.Range("M16").Value = -601
$ws.Range("H31").Value = 5970.9546
$ws.Range("I31").Value = 4329.9
$ws.Range("J31").Value = 7338.5
$ws.Range("K31").Value = 4329.9
$ws.Range("L31").Value = 7338.5
$ws.Range("M31").Value = -4034.9
$ws.Range("N31").Value = -7928.5
$ws.Range("H34").Value = 5970.9546
$ws.Range("I34").Value = 4329.9
$ws.Range("J34").Value = 7338.5
$ws.Range("K34").Value = 4329.9
$ws.Range("L34").Value = 7338.5
$ws.Range("M34").Value = -4127.9
$ws.Range("N34").Value = -7742.5
$ws.Range("H88").Value = 23330
$ws.Range("J88").Value = 23330
$ws.Range("L88").Value = 23330
$ws.Range("N88").Value = -24142
$ws.Range("H91").Value = 23330
$ws.Range("J91").Value = 23330
$ws.Range("L91").Value = 23330
$ws.Range("N91").Value = -26138
$ws.Range("H113").Value = 999.2
$ws.Range("I113").Value = 888
$ws.Range("K113").Value = 888
$ws.Range("M113").Value = 1282
$ws.Range("H141").Value = 269139.1
$ws.Range("J141").Value = 269139.1
$ws.Range("L141").Value = 269139.1
$ws.Range("N141").Value = -279499.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 26.95
$ws.Range("J2").Value = 10.125
$ws.Range("L2").Value = 60.75
$ws.Range("N2").Value = -286.75
$ws.Range("H8").Value = 1597.7273
$ws.Range("I8").Value = 1597.7273
$ws.Range("K8").Value = 4793.1819
$ws.Range("M8").Value = -4654.1819
$ws.Range("H11").Value = 561.875
$ws.Range("I11").Value = 570.7143
$ws.Range("K11").Value = 1712.1429
$ws.Range("M11").Value = -1572.1429
$ws.Range("H38").Value = 1101
$ws.Range("I38").Value = 883.3333
$ws.Range("J38").Value = 1209.8334
$ws.Range("K38").Value = 2649.9999
$ws.Range("L38").Value = 3629.5002
$ws.Range("M38").Value = -2302.9999
$ws.Range("N38").Value = -4323.5002
$ws.Range("H50").Value = 166915
$ws.Range("I50").Value = 290
$ws.Range("K50").Value = 870
$ws.Range("M50").Value = -389
$ws.Range("H53").Value = 166915
$ws.Range("I53").Value = 290
$ws.Range("K53").Value = 870
$ws.Range("M53").Value = -389
$ws.Range("H60").Value = 1678
$ws.Range("I60").Value = 898.25
$ws.Range("K60").Value = 2694.75
$ws.Range("M60").Value = -2443.75
$ws.Range("H98").Value = 1989.3334
$ws.Range("I98").Value = 479
$ws.Range("J98").Value = 2744.5
$ws.Range("K98").Value = 1437
$ws.Range("L98").Value = 8233.5
$ws.Range("M98").Value = 61
$ws.Range("N98").Value = -11229.5
$ws.Range("H107").Value = 553.6923
$ws.Range("J107").Value = 537.3333
$ws.Range("L107").Value = 1611.9999
$ws.Range("N107").Value = -5451.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 33331.332
$ws.Range("J15").Value = 33331.332
$ws.Range("L15").Value = 33331.332
$ws.Range("N15").Value = -33907.332
$ws.Range("H33").Value = 5248.75
$ws.Range("J33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("N33").Value = -3504
$ws.Range("H36").Value = 1694
$ws.Range("I36").Value = 388
$ws.Range("K36").Value = 388
$ws.Range("M36").Value = 97
$ws.Range("H81").Value = 33331.332
$ws.Range("J81").Value = 33331.332
$ws.Range("L81").Value = 33331.332
$ws.Range("N81").Value = -35327.332
$ws.Range("H84").Value = 33331.332
$ws.Range("J84").Value = 33331.332
$ws.Range("L84").Value = 99993.99600000001
$ws.Range("N84").Value = -109977.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1138.8
$ws.Range("I9").Value = 127
$ws.Range("J9").Value = 3499.6667
$ws.Range("K9").Value = 127
$ws.Range("L9").Value = 3499.6667
$ws.Range("M9").Value = 97
$ws.Range("N9").Value = -3947.6667
$ws.Range("H87").Value = 55500
$ws.Range("J87").Value = 55500
$ws.Range("L87").Value = 55500
$ws.Range("N87").Value = -57746
$ws.Range("H90").Value = 55500
$ws.Range("J90").Value = 55500
$ws.Range("L90").Value = 166500
$ws.Range("N90").Value = -177732
$ws.Range("H136").Value = 5724.5
$ws.Range("I136").Value = 3449.5
$ws.Range("K136").Value = 10348.5
$ws.Range("M136").Value = -7798.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 513.2857
$ws.Range("I107").Value = 258.7
$ws.Range("K107").Value = 776.0999999999999
$ws.Range("M107").Value = 1143.9
$ws.Range("H132").Value = 31469.547
$ws.Range("I132").Value = 48135.25
$ws.Range("J132").Value = 9248.611000000001
$ws.Range("K132").Value = 144405.75
$ws.Range("L132").Value = 27745.833
$ws.Range("M132").Value = -141875.75
$ws.Range("N132").Value = -32805.833
$ws.Range("H136").Value = 4289.231
$ws.Range("I136").Value = 3906.238
$ws.Range("K136").Value = 11718.714
$ws.Range("M136").Value = -9168.714

Write-Host "Updated leve profit figures on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
